$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 19 with values for columns A-G
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 25
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 53
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 90
